$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 345 (existing rows 345..366 shift down to 346..367).
$ws.Rows(345).Insert()

# Populate the newly inserted row 345 with the new weekly record.
$ws.Range("A345").Value = 4
$ws.Range("B345").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C345").Value = "Los Lagos"
$ws.Range("D345").Value = 44706
$ws.Range("E345").Value = 10
$ws.Range("F345").Value = 100114013
$ws.Range("G345").Value = "Zanahoria"
$ws.Range("H345").Value = "Sin especificar"
$ws.Range("I345").Value = "Primera"
$ws.Range("J345").Value = 70
$ws.Range("K345").Value = 7000
$ws.Range("L345").Value = 7000
$ws.Range("M345").Value = 7000
$ws.Range("N345").Value = "$/saco 20 kilos"
$ws.Range("O345").Value = "Provincia de Llanquihue"
$ws.Range("P345").Value = 350
$ws.Range("Q345").Value = 20
$ws.Range("R345").Value = "Hortaliza"
